$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
[void]$ws.Activate()

# Build the new "DocumentRepository - Default access" rule row (row 28) by
# cloning the formatting and values of the very similar "Complaint - Default
# access" row (row 20), then overwriting the two cells that differ.
$srcRow = $ws.Range("B20:H20")
$dstRow = $ws.Range("B28:H28")

$srcRow.Copy()
$dstRow.PasteSpecial(-4122)   # xlPasteFormats - bring across styles/borders
$srcRow.Copy()
$dstRow.PasteSpecial(-4163)   # xlPasteValues - bring across the shared values
$excel.CutCopyMode = $false

$ws.Range("B28").Value = "DocumentRepository " + [char]0x2013 + " Default access"
$ws.Range("C28").Value = "DOC_REPO"
# D28 ("participants...isEmpty()") and G28 ("*, *") already match row 20's
# values after the paste above, so nothing else to change there.

# Match the row height used by the other wrapped-text rule rows (30pt).
$ws.Rows.Item(28).RowHeight = 30

# Refresh the view: scroll the window back to the top of the table and move
# the selection down to the newly added row, as happened in the real edit.
[void]$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B29").Select()
